# Apply cryptos list price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.718.16"
$ws.Range("E2").Value = "  -2.80%  "
$ws.Range("D3").Value = "'2.095.13"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'343.76"
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.5163"
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").Value = "'0.4378"
$ws.Range("E8").Value = "  -4.24%  "
$ws.Range("D9").Value = "'52.95"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "'0.09248"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").Value = "'1.163"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "'24.90"
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("D13").Value = "'2.095.38"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "'8.313"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "'6.741"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "'99.38"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'20.80"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "'0.06654"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "'1.009"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'6.184"
$ws.Range("E22").Value = "  -3.14%  "
$ws.Range("D23").Value = "'29.753.56"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("D24").Value = "'12.50"
$ws.Range("E24").Value = "  -3.39%  "
$ws.Range("D25").Value = "'2.319"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("D26").Value = "'2.349.04"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").Value = "'21.95"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").Value = "'2.517"
$ws.Range("D29").Value = "'161.35"
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("D30").Value = "'132.98"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("D31").Value = "'1.131"
$ws.Range("E31").Value = "  -7.67%  "
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("D34").Value = "'6.150"
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("D35").Value = "'3.936"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").Value = "'6.257"
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("D37").Value = "'10.19"
$ws.Range("E37").Value = "  -3.22%  "
$ws.Range("D38").Value = "'0.02575"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").Value = "'0.06694"
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'12.44"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6875"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.326"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "'0.2227"
$ws.Range("E43").Value = "  -4.87%  "
$ws.Range("D44").Value = "'0.6698"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").Value = "'2.316"
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("D47").Value = "'0.00000000359"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").Value = "'3.623"
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("D50").Value = "'81.75"
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'1.161"
$ws.Range("E51").Value = "  -1.94%  "
